$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D2: "0.9000" -> "0.9" (must remain a text value, not be
# re-interpreted as a number by Excel's smart entry) -------------------
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.9"
$ws.Range("A2").Copy()
$ws.Range("D2").PasteSpecial(-4122) # xlPasteFormats - restore original (default) cell style

# --- Add new row 3, reusing row 2's cell formatting so the new cells
# keep the same (default) style as the existing data row ---------------
$ws.Range("A2:O2").Copy()
$ws.Range("A3:O3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "BANKBILL"
$ws.Range("C3").Value = "AUD"

# D3 looks numeric ("0.855"); force it to stay text like the source data
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "0.855"
$ws.Range("A2").Copy()
$ws.Range("D3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "null"

# H3 looks like a date ("2020-02-19"); force it to stay text
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "2020-02-19"
$ws.Range("A2").Copy()
$ws.Range("H3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("I3").Value = "null"
$ws.Range("J3").Value = "null"
$ws.Range("K3").Value = "001M"
$ws.Range("L3").Value = "DEPOSIT"
$ws.Range("M3").Value = "COMRLENDING"
$ws.Range("N3").Value = "null"
$ws.Range("O3").Value = "AUD,EUR"

$excel.CutCopyMode = 0
